# Auto-generated edit script: updates LeveProfits price/profit columns (H-N)
# per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Days of Chunder / Antidote (row 6)
$ws.Range("H6").Value = 788.7
$ws.Range("I6").Value = 800
$ws.Range("K6").Value = 2400
$ws.Range("M6").Value = -2288
# Just Give Him a Serum / Hi-Potion of Strength (row 38)
$ws.Range("H38").Value = 4196.077
$ws.Range("I38").Value = 221.14285
$ws.Range("J38").Value = 8833.5
$ws.Range("K38").Value = 663.4285500000001
$ws.Range("L38").Value = 26500.5
$ws.Range("M38").Value = -291.4285500000001
$ws.Range("N38").Value = -27244.5
# A Bile Business / Shark Oil (row 51)
$ws.Range("H51").Value = 18499.6
$ws.Range("I51").Value = 18499.6
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 18499.6
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -18015.6
$ws.Range("N51").ClearContents()
# The Dotted Line / Enchanted Durium Ink (row 98)
$ws.Range("H98").Value = 1045.3939
$ws.Range("I98").Value = 849.9666999999999
$ws.Range("J98").Value = 2999.6667
$ws.Range("K98").Value = 849.9666999999999
$ws.Range("L98").Value = 2999.6667
$ws.Range("M98").Value = 648.0333000000001
$ws.Range("N98").Value = -5995.6667
# Wishful Inking / Enchanted High Durium Ink (row 122)
$ws.Range("H122").Value = 1045.3939
$ws.Range("I122").Value = 849.9666999999999
$ws.Range("J122").Value = 2999.6667
$ws.Range("K122").Value = 2549.9001
$ws.Range("L122").Value = 8999.000100000001
$ws.Range("M122").Value = -99.90009999999984
$ws.Range("N122").Value = -13899.0001
# Nearly Bare / Gaja Grimoire (row 123)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Fast-forwarding Flora / Growth Formula Lambda (row 132)
$ws.Range("H132").Value = 2534.1943
$ws.Range("I132").Value = 2563.743
$ws.Range("K132").Value = 7691.228999999999
$ws.Range("M132").Value = -5161.228999999999

$ws = $wb.Worksheets.Item("ARM")
# Ain't Got No Ingots / Bronze Ingot (row 2)
$ws.Range("H2").Value = 1696.6875
$ws.Range("I2").Value = 581.4545000000001
$ws.Range("K2").Value = 581.4545000000001
$ws.Range("M2").Value = -468.4545000000001
# No Scope / Titanbronze Ingot (row 116)
$ws.Range("H116").Value = 1696.6875
$ws.Range("I116").Value = 581.4545000000001
$ws.Range("K116").Value = 581.4545000000001
$ws.Range("M116").Value = 1712.5455
# Don't Bore Me, Ore Me / Mountain Chromite Ingot (row 132)
$ws.Range("H132").Value = 3671.353
$ws.Range("I132").Value = 2022.7407
$ws.Range("J132").Value = 10030.286
$ws.Range("K132").Value = 6068.2221
$ws.Range("L132").Value = 30090.858
$ws.Range("M132").Value = -3538.2221
$ws.Range("N132").Value = -35150.858

$ws = $wb.Worksheets.Item("BSM")
# Hells Bells / Bronze Ingot (row 3)
$ws.Range("H3").Value = 1696.6875
$ws.Range("I3").Value = 581.4545000000001
$ws.Range("K3").Value = 581.4545000000001
$ws.Range("M3").Value = -467.4545000000001
# Through Thick and Thin / Adamantite Nugget (row 86)
$ws.Range("H86").Value = 4276.55
$ws.Range("I86").Value = 2404.5833
$ws.Range("K86").Value = 2404.5833
$ws.Range("M86").Value = -1281.5833
# Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget (row 89)
$ws.Range("H89").Value = 4276.55
$ws.Range("I89").Value = 2404.5833
$ws.Range("K89").Value = 12022.9165
$ws.Range("M89").Value = -6406.916499999999
# Meddle in Metal / Oroshigane Ingot (row 99)
$ws.Range("H99").Value = 2364.9167
$ws.Range("I99").Value = 2533.5715
$ws.Range("K99").Value = 2533.5715
$ws.Range("M99").Value = -1035.5715
# Ruthenium Supremium / Ruthenium Ingot (row 134)
$ws.Range("H134").Value = 7080.294
$ws.Range("I134").Value = 4613.864
$ws.Range("J134").Value = 11602.083
$ws.Range("K134").Value = 13841.592
$ws.Range("L134").Value = 34806.249
$ws.Range("M134").Value = -11306.592
$ws.Range("N134").Value = -39876.249

$ws = $wb.Worksheets.Item("CRP")
# Compulsory Conjury / Maple Cane (row 13)
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# O Pine / Pine Lumber (row 99)
$ws.Range("H99").Value = 5381.75
$ws.Range("I99").Value = 7389.476
$ws.Range("J99").Value = 3548.6086
$ws.Range("K99").Value = 7389.476
$ws.Range("L99").Value = 3548.6086
$ws.Range("M99").Value = -5891.476
$ws.Range("N99").Value = -6544.6086
# Timber of Tenkonto / Horse Chestnut Lumber (row 122)
$ws.Range("H122").Value = 3176.7036
$ws.Range("J122").Value = 4510.6
$ws.Range("L122").Value = 13531.8
$ws.Range("N122").Value = -18431.8
# A Better Conductor / Red Pine Lumber (row 126)
$ws.Range("H126").Value = 5381.75
$ws.Range("I126").Value = 7389.476
$ws.Range("J126").Value = 3548.6086
$ws.Range("K126").Value = 22168.428
$ws.Range("L126").Value = 10645.8258
$ws.Range("M126").Value = -19698.428
$ws.Range("N126").Value = -15585.8258
# Hull Lotta Damage / Ginseng Lumber (row 132)
$ws.Range("H132").Value = 3292.7292
$ws.Range("I132").Value = 3126.7673
$ws.Range("J132").Value = 4720
$ws.Range("K132").Value = 9380.3019
$ws.Range("L132").Value = 14160
$ws.Range("M132").Value = -6850.3019
$ws.Range("N132").Value = -19220
# Wood You Be Quiet / Ceiba Lumber (row 134)
$ws.Range("H134").Value = 2487.5122
$ws.Range("I134").Value = 2035.5312
$ws.Range("J134").Value = 4094.5557
$ws.Range("K134").Value = 6106.5936
$ws.Range("L134").Value = 12283.6671
$ws.Range("M134").Value = -3571.5936
$ws.Range("N134").Value = -17353.6671

$ws = $wb.Worksheets.Item("CUL")
# Salt of the North / Northern Sea Salt (row 122)
$ws.Range("H122").Value = 502
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 502
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4518
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9418

$ws = $wb.Worksheets.Item("GSM")
# Needful Rings / Copper Wristlets (row 3)
$ws.Range("H3").Value = 61341.8
$ws.Range("I3").Value = 38333
$ws.Range("K3").Value = 38333
$ws.Range("M3").Value = -38217
# Needs More Prayerbell / Hardsilver Ingot (row 80)
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# With a Noise That Reaches Heaven (L) / Hardsilver Ingot (row 83)
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# If I'd a Koppranickel for Every Time... / Koppranickel Ingot (row 97)
$ws.Range("H97").Value = 834.58826
$ws.Range("I97").Value = 823.5454999999999
$ws.Range("J97").Value = 854.8333
$ws.Range("K97").Value = 823.5454999999999
$ws.Range("L97").Value = 854.8333
$ws.Range("M97").Value = -327.5454999999999
$ws.Range("N97").Value = -1846.8333
# Awarding Academic Excellence / Ametrine (row 122)
$ws.Range("H122").Value = 3277.25
$ws.Range("I122").Value = 2223.3901
$ws.Range("J122").Value = 5551.3687
$ws.Range("K122").Value = 6670.1703
$ws.Range("L122").Value = 16654.1061
$ws.Range("M122").Value = -4220.1703
$ws.Range("N122").Value = -21554.1061
# Gold Rush Order / Phrygian Gold Ingot (row 126)
$ws.Range("H126").Value = 6515.231
$ws.Range("I126").Value = 6124.5
$ws.Range("K126").Value = 18373.5
$ws.Range("M126").Value = -15903.5

$ws = $wb.Worksheets.Item("LTW")
# Best Served Toad / Toad Leather (row 40)
$ws.Range("H40").Value = 7467
$ws.Range("I40").Value = 3826
$ws.Range("J40").Value = 14749
$ws.Range("K40").Value = 3826
$ws.Range("L40").Value = 14749
$ws.Range("M40").Value = -3690
$ws.Range("N40").Value = -15021
# Trainin' the Neck / Dragon Leather (row 82)
$ws.Range("H82").Value = 3117
$ws.Range("I82").Value = 1320.5
$ws.Range("J82").Value = 4314.6665
$ws.Range("K82").Value = 1320.5
$ws.Range("L82").Value = 4314.6665
$ws.Range("M82").Value = -959.5
$ws.Range("N82").Value = -5036.6665
# Training Is Only Skintight (L) / Dragon Leather (row 85)
$ws.Range("H85").Value = 3117
$ws.Range("I85").Value = 1320.5
$ws.Range("J85").Value = 4314.6665
$ws.Range("K85").Value = 1320.5
$ws.Range("L85").Value = 4314.6665
$ws.Range("M85").Value = -72.5
$ws.Range("N85").Value = -6810.6665
# Hell on Leather / Gaja Leather (row 122)
$ws.Range("H122").Value = 3782.7896
$ws.Range("I122").Value = 2648.625
$ws.Range("J122").Value = 9831.666999999999
$ws.Range("K122").Value = 7945.875
$ws.Range("L122").Value = 29495.001
$ws.Range("M122").Value = -5495.875
$ws.Range("N122").Value = -34395.001

$ws = $wb.Worksheets.Item("WVR")
# Comfy Cabins / Snow Cotton Cloth (row 132)
$ws.Range("H132").Value = 2749.5715
$ws.Range("I132").Value = 2551.1133
$ws.Range("K132").Value = 7653.3399
$ws.Range("M132").Value = -5123.3399
